$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new row 2 for "2022-Q3" and push
#    the existing quarters down by one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

# restore the index-column style (bold/centered/bordered) on A2 by
# copying it from A3 (an existing, already-shifted data row)
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 0.54

# the "index" column (A) is a 0-based running count; bump every
# pushed-down row by one so it keeps counting from the new top row
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(7,1).Value = 5

# ------------------------------------------------------------------
# 2) Add the new "2022-Q3" detail sheet, right before "2022-Q2".
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q2"))
$newSheet.Name = "2022-Q3"

# Copy header/index-column style (bold/centered/bordered) from the
# "2022-Q2" sheet so the new sheet matches the workbook's look. Fetch
# the source sheet fresh (post-insert) since earlier handles go stale
# once the workbook's sheet collection changes.
$q2ws = $wb.Worksheets.Item("2022-Q2")
$q2ws.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q2ws.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows (leading apostrophe forces text storage for numeric-looking values)
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'960004"
$newSheet.Cells.Item(2,3).Value = "华夏兴华混合H"
$newSheet.Cells.Item(2,4).Value = "'10.35"
$newSheet.Cells.Item(2,5).Value = "'83.78"
$newSheet.Cells.Item(2,6).Value = "'3.35"
$newSheet.Cells.Item(2,7).Value = "'0.3467"
$newSheet.Cells.Item(2,8).Value = 9

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'006868"
$newSheet.Cells.Item(3,3).Value = "华夏科技成长股票"
$newSheet.Cells.Item(3,4).Value = "'5.37"
$newSheet.Cells.Item(3,5).Value = "'88.33"
$newSheet.Cells.Item(3,6).Value = "'3.51"
$newSheet.Cells.Item(3,7).Value = "'0.1885"
$newSheet.Cells.Item(3,8).Value = 8

$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "'519908"
$newSheet.Cells.Item(4,3).Value = "华夏兴华混合A"
$newSheet.Cells.Item(4,4).Value = "'0.00"
$newSheet.Cells.Item(4,5).Value = "'83.78"
$newSheet.Cells.Item(4,6).Value = "'3.35"
$newSheet.Cells.Item(4,7).Value = 0
$newSheet.Cells.Item(4,8).Value = 9

# The leading apostrophes above force text storage for numeric-looking
# values, but they also leave a "number stored as text" (quote-prefix)
# cell style behind. Strip that cosmetic style back off (this keeps
# the text type/value, only the formatting resets to plain/default) so
# the data rows stay unstyled, matching the rest of the workbook.
$newSheet.Range("B2:H4").ClearFormats()
